$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'20.080.45"
$ws.Range('E2').Value = '  -7.79%  '
$ws.Range('D3').Value = "'1.428.30"
$ws.Range('E3').Value = '  -7.38%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').Value = "'274.58"
$ws.Range('E6').Value = '  -5.25%  '
$ws.Range('D7').Value = "'0.3742"
$ws.Range('E7').Value = '  -3.75%  '
$ws.Range('D8').Value = "'0.3090"
$ws.Range('E8').Value = '  -3.28%  '
$ws.Range('D9').Value = "'40.13"
$ws.Range('E9').Value = '  -7.64%  '
$ws.Range('D10').Value = "'1.016"
$ws.Range('E10').Value = '  -4.62%  '
$ws.Range('D11').Value = "'0.06588"
$ws.Range('E11').Value = '  -8.45%  '
$ws.Range('D12').Value = "'1.002"
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').Value = "'5.387"
$ws.Range('E13').Value = '  -4.43%  '
$ws.Range('D14').Value = "'17.18"
$ws.Range('E14').Value = '  -7.57%  '
$ws.Range('D15').Value = "'6.175"
$ws.Range('E15').Value = '  -6.56%  '
$ws.Range('D16').Value = "'1.430.48"
$ws.Range('E16').Value = '  -7.20%  '
$ws.Range('D17').Value = "'0.00001010"
$ws.Range('E17').Value = '  -8.69%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = "'0.05805"
$ws.Range('E18').Value = '  -11.87%  '
$ws.Range('B19').Value = 'Litecoin'
$ws.Range('C19').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D19').Value = "'75.40"
$ws.Range('E19').Value = '  -9.56%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').Value = "'5.684"
$ws.Range('E21').Value = '  -7.44%  '
$ws.Range('D22').Value = "'14.51"
$ws.Range('E22').Value = '  -5.72%  '
$ws.Range('D23').Value = "'11.09"
$ws.Range('E23').Value = '  +1.98%  '
$ws.Range('D24').Value = "'2.338"
$ws.Range('E24').Value = '  -1.22%  '
$ws.Range('D25').Value = "'20.101.52"
$ws.Range('E25').Value = '  -7.76%  '
$ws.Range('D26').Value = "'2.276"
$ws.Range('E26').Value = '  -4.84%  '
$ws.Range('D27').Value = "'138.16"
$ws.Range('E27').Value = '  -5.02%  '
$ws.Range('D28').Value = "'16.89"
$ws.Range('E28').Value = '  -8.36%  '
$ws.Range('D29').Value = "'1.591.10"
$ws.Range('E29').Value = '  -7.32%  '
$ws.Range('D30').Value = "'109.27"
$ws.Range('E30').Value = '  -7.32%  '
$ws.Range('D31').Value = "'3.957"
$ws.Range('E31').Value = '  -18.47%  '
$ws.Range('D32').Value = "'0.9044"
$ws.Range('E32').Value = '  -6.21%  '
$ws.Range('D33').Value = "'5.409"
$ws.Range('E33').Value = '  -8.33%  '
$ws.Range('D34').Value = "'0.07783"
$ws.Range('E34').Value = '  -5.04%  '
$ws.Range('D35').Value = "'8.387"
$ws.Range('E35').Value = '  -6.51%  '
$ws.Range('D36').Value = "'11.37"
$ws.Range('E36').Value = '  +6.36%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').Value = "'1.001"
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = "'4.763"
$ws.Range('E38').Value = '  -7.16%  '
$ws.Range('D39').Value = "'0.05700"
$ws.Range('E39').Value = '  -7.04%  '
$ws.Range('D40').Value = "'0.1921"
$ws.Range('E40').Value = '  -6.05%  '
$ws.Range('D41').Value = "'1.116"
$ws.Range('E41').Value = '  -5.66%  '
$ws.Range('D42').Value = "'0.02028"
$ws.Range('E42').Value = '  -8.34%  '
$ws.Range('D43').Value = "'1.307"
$ws.Range('E43').Value = '  -10.13%  '
$ws.Range('D44').Value = "'0.5331"
$ws.Range('E44').Value = '  -7.50%  '
$ws.Range('D45').Value = "'3.543"
$ws.Range('E45').Value = '  -5.28%  '
$ws.Range('D46').Value = "'12.14"
$ws.Range('E46').Value = '  -7.20%  '
$ws.Range('D47').Value = "'0.5140"
$ws.Range('E47').Value = '  -7.09%  '
$ws.Range('D48').Value = "'1.777"
$ws.Range('E48').Value = '  -5.32%  '
$ws.Range('D49').Value = "'109.47"
$ws.Range('E49').Value = '  -7.28%  '
$ws.Range('D50').Value = "'1.051"
$ws.Range('E50').Value = '  -7.42%  '
$ws.Range('D51').Value = "'1.001"
$ws.Range('E51').Value = '  +0.10%  '